$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.407.31'
$ws.Range("E2").Value = '  -0.48%  '

$ws.Range("D3").Value = '1.847.10'
$ws.Range("E3").Value = '  -0.30%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9984'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.30'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.17%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6332'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.07%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.0000'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("E8").Value = '  -0.16%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2966'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.49%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.63'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.23%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07732'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.61%  '

$ws.Range("D12").Value = '1.855.37'
$ws.Range("E12").Value = '  -0.33%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.001'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.70%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6851'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.42%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '83.07'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.05%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000009973'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.10%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.179'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.99%  '

$ws.Range("D18").Value = '29.424.88'
$ws.Range("E18").Value = '  -0.57%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '230.05'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.26%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.48'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.39%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9999'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.01%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.574'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.02%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.959'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.66%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.9999'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.01%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '157.14'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.76%  '

$ws.Range("E26").Value = '  +1.01%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.384'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.05%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.68'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.47%  '

$ws.Range("E29").Value = '  -1.36%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05729'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.94%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.252'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.53%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.131'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.26%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.037'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.42%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.854'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.38%  '

$ws.Range("E35").Value = '  -1.30%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7167'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.41%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.595'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.03%  '

$ws.Range("D38").Value = '1.252.19'
$ws.Range("E38").Value = '  +1.14%  '

$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01810'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.86%  '

$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.782'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.76%  '

$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.210'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.07%  '

$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9086'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.75%  '

$ws.Range("E43").Value = '  +0.10%  '

$ws.Range("D44").Value = '2.011.79'
$ws.Range("E44").Value = '  -1.35%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '101.80'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.19%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '66.56'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.72%  '

$ws.Range("B47").Value = 'Aptos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.093'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.87%  '

$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.166'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.24%  '

$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.714'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.81%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4030'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.44%  '

$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.00000000116'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.79%  '
